$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cycle timing labels in row 6: the "A: Day 1", "A: Day 16" and
# "P4: +30 Days" timings are all replaced with the generic cycle marker "C:"
$ws.Range("F6").Value = "C:"
$ws.Range("H6").Value = "C:"
$ws.Range("J6").Value = "C:"

# Update the active selection to H6, matching the saved view state
$ws.Range("H6").Select()
